$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 87.28570999999999
$ws.Cells.Item(6, 9).Value = 18.5
$ws.Cells.Item(6, 11).Value = 55.5
$ws.Cells.Item(6, 13).Value = 56.5
$ws.Cells.Item(8, 8).Value = 957
$ws.Cells.Item(8, 9).Value = 103.6
$ws.Cells.Item(8, 10).Value = 1810.4
$ws.Cells.Item(8, 11).Value = 310.8
$ws.Cells.Item(8, 12).Value = 5431.200000000001
$ws.Cells.Item(8, 13).Value = -171.8
$ws.Cells.Item(8, 14).Value = -5709.200000000001
$ws.Cells.Item(12, 8).Value = 839.5
$ws.Cells.Item(12, 9).Value = 785.8333
$ws.Cells.Item(12, 11).Value = 785.8333
$ws.Cells.Item(12, 13).Value = -615.8333
$ws.Cells.Item(17, 8).Value = 1646.1111
$ws.Cells.Item(17, 10).Value = 1646.1111
$ws.Cells.Item(17, 12).Value = 4938.3333
$ws.Cells.Item(17, 14).Value = -5274.3333
$ws.Cells.Item(33, 8).Value = 756
$ws.Cells.Item(33, 9).Value = 1034
$ws.Cells.Item(33, 10).Value = 200
$ws.Cells.Item(33, 11).Value = 1034
$ws.Cells.Item(33, 12).Value = 200
$ws.Cells.Item(33, 13).Value = -805
$ws.Cells.Item(33, 14).Value = -658
$ws.Cells.Item(40, 8).Value = 7731.636
$ws.Cells.Item(40, 9).Value = 5925.3335
$ws.Cells.Item(40, 11).Value = 5925.3335
$ws.Cells.Item(40, 13).Value = -5750.3335
$ws.Cells.Item(93, 8).Value = 28333
$ws.Cells.Item(93, 10).Value = 28333
$ws.Cells.Item(93, 12).Value = 28333
$ws.Cells.Item(93, 14).Value = -33325
$ws.Cells.Item(98, 8).Value = 989.3
$ws.Cells.Item(98, 9).Value = 642
$ws.Cells.Item(98, 11).Value = 642
$ws.Cells.Item(98, 13).Value = 856
$ws.Cells.Item(122, 8).Value = 989.3
$ws.Cells.Item(122, 9).Value = 642
$ws.Cells.Item(122, 11).Value = 1926
$ws.Cells.Item(122, 13).Value = 524
$ws.Cells.Item(132, 8).Value = 3996.9333
$ws.Cells.Item(132, 9).Value = 1329.5
$ws.Cells.Item(132, 11).Value = 3988.5
$ws.Cells.Item(132, 13).Value = -1458.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 6879.8887
$ws.Cells.Item(2, 9).Value = 1070.1666
$ws.Cells.Item(2, 11).Value = 1070.1666
$ws.Cells.Item(2, 13).Value = -957.1666
$ws.Cells.Item(44, 8).Value = 11861.277
$ws.Cells.Item(44, 10).Value = 12441.353
$ws.Cells.Item(44, 12).Value = 12441.353
$ws.Cells.Item(44, 14).Value = -13417.353
$ws.Cells.Item(45, 8).Value = 3415.0715
$ws.Cells.Item(45, 9).Value = 2444.1667
$ws.Cells.Item(45, 11).Value = 2444.1667
$ws.Cells.Item(45, 13).Value = -2067.1667
$ws.Cells.Item(61, 8).Value = 1797.7059
$ws.Cells.Item(61, 10).Value = 3999.5
$ws.Cells.Item(61, 12).Value = 3999.5
$ws.Cells.Item(61, 14).Value = -4423.5
$ws.Cells.Item(102, 8).Value = 11367609
$ws.Cells.Item(102, 9).Value = 25001188
$ws.Cells.Item(102, 10).Value = 6291.6665
$ws.Cells.Item(102, 11).Value = 25001188
$ws.Cells.Item(102, 12).Value = 6291.6665
$ws.Cells.Item(102, 13).Value = -24999566
$ws.Cells.Item(102, 14).Value = -9535.666499999999
$ws.Cells.Item(116, 8).Value = 6879.8887
$ws.Cells.Item(116, 9).Value = 1070.1666
$ws.Cells.Item(116, 11).Value = 1070.1666
$ws.Cells.Item(116, 13).Value = 1223.8334
$ws.Cells.Item(132, 8).Value = 1121.125
$ws.Cells.Item(132, 9).Value = 1121.125
$ws.Cells.Item(132, 11).Value = 3363.375
$ws.Cells.Item(132, 13).Value = -833.375
$ws.Cells.Item(136, 8).Value = 1797.7059
$ws.Cells.Item(136, 10).Value = 3999.5
$ws.Cells.Item(136, 12).Value = 11998.5
$ws.Cells.Item(136, 14).Value = -17098.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 6879.8887
$ws.Cells.Item(3, 9).Value = 1070.1666
$ws.Cells.Item(3, 11).Value = 1070.1666
$ws.Cells.Item(3, 13).Value = -956.1666
$ws.Cells.Item(5, 8).Value = 168
$ws.Cells.Item(5, 9).Value = 102
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 102
$ws.Cells.Item(5, 12).Value = 300
$ws.Cells.Item(5, 13).Value = 11
$ws.Cells.Item(5, 14).Value = -526
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 54.153847
$ws.Cells.Item(7, 9).Value = 53.11111
$ws.Cells.Item(7, 10).Value = 56.5
$ws.Cells.Item(7, 11).Value = 53.11111
$ws.Cells.Item(7, 12).Value = 56.5
$ws.Cells.Item(7, 13).Value = 59.88889
$ws.Cells.Item(7, 14).Value = -282.5
$ws.Cells.Item(13, 8).Value = 476.16666
$ws.Cells.Item(13, 10).Value = 552.3333
$ws.Cells.Item(13, 12).Value = 552.3333
$ws.Cells.Item(13, 14).Value = -830.3333
$ws.Cells.Item(22, 8).Value = 775.3158
$ws.Cells.Item(22, 9).Value = 649.5333000000001
$ws.Cells.Item(22, 10).Value = 1247
$ws.Cells.Item(22, 11).Value = 649.5333000000001
$ws.Cells.Item(22, 12).Value = 1247
$ws.Cells.Item(22, 13).Value = -299.5333000000001
$ws.Cells.Item(22, 14).Value = -1947
$ws.Cells.Item(58, 8).Value = 3903.7693
$ws.Cells.Item(58, 9).Value = 1727.75
$ws.Cells.Item(58, 10).Value = 7385.4
$ws.Cells.Item(58, 11).Value = 1727.75
$ws.Cells.Item(58, 12).Value = 7385.4
$ws.Cells.Item(58, 13).Value = -1524.75
$ws.Cells.Item(58, 14).Value = -7791.4
$ws.Cells.Item(94, 8).Value = 14507
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 14507
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 14507
$ws.Cells.Item(94, 13).Value = $null
$ws.Cells.Item(94, 14).Value = -15409
$ws.Cells.Item(99, 8).Value = 2640
$ws.Cells.Item(99, 9).Value = 2174.6
$ws.Cells.Item(99, 11).Value = 2174.6
$ws.Cells.Item(99, 13).Value = -676.5999999999999
$ws.Cells.Item(122, 8).Value = 1332
$ws.Cells.Item(122, 9).Value = 1516.6666
$ws.Cells.Item(122, 10).Value = 962.6667
$ws.Cells.Item(122, 11).Value = 4549.9998
$ws.Cells.Item(122, 12).Value = 2888.0001
$ws.Cells.Item(122, 13).Value = -2099.9998
$ws.Cells.Item(122, 14).Value = -7788.0001
$ws.Cells.Item(126, 8).Value = 2640
$ws.Cells.Item(126, 9).Value = 2174.6
$ws.Cells.Item(126, 11).Value = 6523.799999999999
$ws.Cells.Item(126, 13).Value = -4053.799999999999
$ws.Cells.Item(134, 8).Value = 3244
$ws.Cells.Item(136, 8).Value = 3903.7693
$ws.Cells.Item(136, 9).Value = 1727.75
$ws.Cells.Item(136, 10).Value = 7385.4
$ws.Cells.Item(136, 11).Value = 5183.25
$ws.Cells.Item(136, 12).Value = 22156.2
$ws.Cells.Item(136, 13).Value = -2633.25
$ws.Cells.Item(136, 14).Value = -27256.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 120.35294
$ws.Cells.Item(2, 9).Value = 149
$ws.Cells.Item(2, 10).Value = 88.125
$ws.Cells.Item(2, 11).Value = 894
$ws.Cells.Item(2, 12).Value = 528.75
$ws.Cells.Item(2, 13).Value = -781
$ws.Cells.Item(2, 14).Value = -754.75
$ws.Cells.Item(34, 8).Value = 1548.3572
$ws.Cells.Item(34, 10).Value = 2915
$ws.Cells.Item(34, 12).Value = 8745
$ws.Cells.Item(34, 14).Value = -8913
$ws.Cells.Item(120, 8).Value = 1499
$ws.Cells.Item(120, 9).Value = 1499
$ws.Cells.Item(120, 11).Value = 4497
$ws.Cells.Item(120, 13).Value = 341
$ws.Cells.Item(122, 8).Value = 799.6667
$ws.Cells.Item(122, 9).Value = 456.7143
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 4110.428699999999
$ws.Cells.Item(122, 12).Value = 18000
$ws.Cells.Item(122, 13).Value = -1660.428699999999
$ws.Cells.Item(122, 14).Value = -22900
$ws.Cells.Item(133, 8).Value = 4000
$ws.Cells.Item(133, 9).Value = 3000
$ws.Cells.Item(133, 11).Value = 9000
$ws.Cells.Item(133, 13).Value = -3940
$ws.Cells.Item(136, 8).Value = 6725
$ws.Cells.Item(136, 9).Value = 450
$ws.Cells.Item(136, 11).Value = 1350
$ws.Cells.Item(136, 13).Value = 3750
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 522.36365
$ws.Cells.Item(2, 9).Value = 249.42857
$ws.Cells.Item(2, 10).Value = 1000
$ws.Cells.Item(2, 11).Value = 249.42857
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = -136.42857
$ws.Cells.Item(2, 14).Value = -1226
$ws.Cells.Item(102, 8).Value = 3631.8823
$ws.Cells.Item(102, 9).Value = 3285.818
$ws.Cells.Item(102, 11).Value = 3285.818
$ws.Cells.Item(102, 13).Value = -1663.818
$ws.Cells.Item(113, 8).Value = 7257
$ws.Cells.Item(113, 10).Value = 10000
$ws.Cells.Item(113, 12).Value = 10000
$ws.Cells.Item(113, 14).Value = -14340
$ws.Cells.Item(122, 8).Value = 3993.7144
$ws.Cells.Item(122, 9).Value = 3997.6
$ws.Cells.Item(122, 10).Value = 3984
$ws.Cells.Item(122, 11).Value = 11992.8
$ws.Cells.Item(122, 12).Value = 11952
$ws.Cells.Item(122, 13).Value = -9542.799999999999
$ws.Cells.Item(122, 14).Value = -16852
$ws.Cells.Item(132, 8).Value = 78378.92
$ws.Cells.Item(132, 9).Value = 84836
$ws.Cells.Item(132, 11).Value = 254508
$ws.Cells.Item(132, 13).Value = -251978
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1139.5
$ws.Cells.Item(22, 9).Value = 832.6667
$ws.Cells.Item(22, 11).Value = 832.6667
$ws.Cells.Item(22, 13).Value = -537.6667
$ws.Cells.Item(27, 8).Value = 1139.5
$ws.Cells.Item(27, 9).Value = 832.6667
$ws.Cells.Item(27, 11).Value = 832.6667
$ws.Cells.Item(27, 13).Value = -725.6667
$ws.Cells.Item(93, 8).Value = 1316.2106
$ws.Cells.Item(93, 9).Value = 1450.6666
$ws.Cells.Item(93, 10).Value = 1085.7142
$ws.Cells.Item(93, 11).Value = 1450.6666
$ws.Cells.Item(93, 12).Value = 1085.7142
$ws.Cells.Item(93, 13).Value = -202.6666
$ws.Cells.Item(93, 14).Value = -3581.7142
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1916.3334
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 1916.3334
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 5749.0002
$ws.Cells.Item(122, 13).Value = $null
$ws.Cells.Item(122, 14).Value = -10649.0002
$ws.Cells.Item(126, 8).Value = 3196.3572
$ws.Cells.Item(126, 9).Value = 527.6667
$ws.Cells.Item(126, 10).Value = 8000
$ws.Cells.Item(126, 11).Value = 1583.0001
$ws.Cells.Item(126, 12).Value = 24000
$ws.Cells.Item(126, 13).Value = 886.9999
$ws.Cells.Item(126, 14).Value = -28940
